$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the D, M, N, O, P, R, S values between row pairs (2,6), (3,7), (4,8), (5,9)
$pairs = @(
    @(2, 6),
    @(3, 7),
    @(4, 8),
    @(5, 9)
)

$cols = @("D", "M", "N", "O", "P", "R", "S")

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"

        $v1 = $ws.Range($addr1).Value()
        $v2 = $ws.Range($addr2).Value()

        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }
}
